$wb = $excel.ActiveWorkbook

# --- 1. Insert a new worksheet "2022-Q4" right before the existing "2022-Q3" sheet ---
# Final sheet order must become: 总计, 2022-Q4, 2022-Q3, 2021-Q2
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q4Sheet = $wb.Worksheets.Add($q3Sheet)
$q4Sheet.Name = "2022-Q4"

# --- 2. Populate the new "2022-Q4" sheet with its fund-holdings data ---

# Header row
$q4Sheet.Range("B1").Value = "基金代码"
$q4Sheet.Range("C1").Value = "基金名称"
$q4Sheet.Range("D1").Value = "基金规模"
$q4Sheet.Range("E1").Value = "股票总仓位"
$q4Sheet.Range("F1").Value = "仓位占比"
$q4Sheet.Range("G1").Value = "持有市值(亿元)"
$q4Sheet.Range("H1").Value = "仓位排名"

# Columns B..G are stored as text in the source data (fund codes have leading
# zeros, figures keep trailing zeros like "28.90") - force text formatting
# before writing so those are preserved instead of being coerced to numbers.
$q4Sheet.Range("B2:G7").NumberFormat = "@"

$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("B2").Value = "630010"
$q4Sheet.Range("C2").Value = "华商价值精选混合"
$q4Sheet.Range("D2").Value = "4.37"
$q4Sheet.Range("E2").Value = "87.25"
$q4Sheet.Range("F2").Value = "3.90"
$q4Sheet.Range("G2").Value = "0.1704"
$q4Sheet.Range("H2").Value = 5

$q4Sheet.Range("A3").Value = 1
$q4Sheet.Range("B3").Value = "630006"
$q4Sheet.Range("C3").Value = "华商产业升级混合"
$q4Sheet.Range("D3").Value = "0.86"
$q4Sheet.Range("E3").Value = "88.65"
$q4Sheet.Range("F3").Value = "3.97"
$q4Sheet.Range("G3").Value = "0.0341"
$q4Sheet.Range("H3").Value = 5

$q4Sheet.Range("A4").Value = 2
$q4Sheet.Range("B4").Value = "011685"
$q4Sheet.Range("C4").Value = "创金合信先进装备股票A"
$q4Sheet.Range("D4").Value = "0.24"
$q4Sheet.Range("E4").Value = "80.29"
$q4Sheet.Range("F4").Value = "9.34"
$q4Sheet.Range("G4").Value = "0.0224"
$q4Sheet.Range("H4").Value = 4

$q4Sheet.Range("A5").Value = 3
$q4Sheet.Range("B5").Value = "011686"
$q4Sheet.Range("C5").Value = "创金合信先进装备股票C"
$q4Sheet.Range("D5").Value = "0.18"
$q4Sheet.Range("E5").Value = "80.29"
$q4Sheet.Range("F5").Value = "9.34"
$q4Sheet.Range("G5").Value = "0.0168"
$q4Sheet.Range("H5").Value = 4

$q4Sheet.Range("A6").Value = 4
$q4Sheet.Range("B6").Value = "012287"
$q4Sheet.Range("C6").Value = "东海启航6个月持有期混合A"
$q4Sheet.Range("D6").Value = "0.74"
$q4Sheet.Range("E6").Value = "28.90"
$q4Sheet.Range("F6").Value = "1.02"
$q4Sheet.Range("G6").Value = "0.0075"
$q4Sheet.Range("H6").Value = 1

$q4Sheet.Range("A7").Value = 5
$q4Sheet.Range("B7").Value = "013377"
$q4Sheet.Range("C7").Value = "东海启航6个月持有期混合C"
$q4Sheet.Range("D7").Value = "0.23"
$q4Sheet.Range("E7").Value = "28.90"
$q4Sheet.Range("F7").Value = "1.02"
$q4Sheet.Range("G7").Value = "0.0023"
$q4Sheet.Range("H7").Value = 1

# Style the header row and the index column like the other quarter sheets
# (bold font, thin border, centered) so the new sheet matches the workbook's look.
$q4Sheet.Range("B1:H1").Font.Bold = $true
$q4Sheet.Range("B1:H1").HorizontalAlignment = -4108
$q4Sheet.Range("B1:H1").VerticalAlignment = -4160
$q4Sheet.Range("B1:H1").Borders.LineStyle = 1

$q4Sheet.Range("A2:A7").Font.Bold = $true
$q4Sheet.Range("A2:A7").HorizontalAlignment = -4108
$q4Sheet.Range("A2:A7").VerticalAlignment = -4160
$q4Sheet.Range("A2:A7").Borders.LineStyle = 1

# --- 3. Update the "总计" (summary) sheet: add the 2022-Q4 row at the top of the
#        data, pushing the existing 2022-Q3 / 2021-Q2 rows down by one ---
$totalSheet = $wb.Worksheets.Item("总计")

# Read the existing values first (Value2, not Value - Value getter is unreliable here)
$r2b = $totalSheet.Range("B2").Value2
$r2c = $totalSheet.Range("C2").Value2
$r2d = $totalSheet.Range("D2").Value2

$r3b = $totalSheet.Range("B3").Value2
$r3c = $totalSheet.Range("C3").Value2
$r3d = $totalSheet.Range("D3").Value2

# Give the new row 4 the same column-A style as the existing index cells
$totalSheet.Range("A3").Copy($totalSheet.Range("A4"))

# Row 4 = old row 3 (2021-Q2)
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = $r3b
$totalSheet.Range("C4").Value = $r3c
$totalSheet.Range("D4").Value = $r3d

# Row 3 = old row 2 (2022-Q3)
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = $r2b
$totalSheet.Range("C3").Value = $r2c
$totalSheet.Range("D3").Value = $r2d

# Row 2 = new (2022-Q4)
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 0.25
